# Auto-generated Excel COM-interop script
# Applies betting odds updates for 2025-10-15 workbook:
#  1) Updates existing cell values in rows 2-7
#  2) Appends 4 new match rows (8-11)
#  3) Dimension auto-extends to A1:AO11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update existing values in rows 2-7 ---
$ws.Cells.Item(2, 6).Value = 1.61
$ws.Cells.Item(2, 7).Value = 1.99
$ws.Cells.Item(2, 8).Value = 4.6
$ws.Cells.Item(2, 9).Value = 7.2
$ws.Cells.Item(2, 10).Value = 3.3
$ws.Cells.Item(2, 11).Value = 6.6
$ws.Cells.Item(2, 12).Value = 1.33
$ws.Cells.Item(2, 13).Value = 1.06
$ws.Cells.Item(2, 14).Value = 3.05
$ws.Cells.Item(2, 15).Value = 1.32
$ws.Cells.Item(2, 16).Value = 1.72
$ws.Cells.Item(2, 17).Value = 1.82
$ws.Cells.Item(2, 18).Value = 1.32
$ws.Cells.Item(2, 19).Value = 3.05
$ws.Cells.Item(2, 20).Value = 1.89
$ws.Cells.Item(2, 21).Value = 1.9
$ws.Cells.Item(2, 22).Value = 1.16
$ws.Cells.Item(2, 23).Value = 2.12
$ws.Cells.Item(2, 28).Value = 9.800000000000001
$ws.Cells.Item(3, 6).Value = 1.36
$ws.Cells.Item(3, 19).Value = 3.15
$ws.Cells.Item(4, 8).Value = 2.1
$ws.Cells.Item(4, 9).Value = 2.26
$ws.Cells.Item(4, 15).Value = 1.44
$ws.Cells.Item(4, 37).Value = 65
$ws.Cells.Item(5, 17).Value = 2.12
$ws.Cells.Item(5, 19).Value = 3.95
$ws.Cells.Item(5, 21).Value = 1.99
$ws.Cells.Item(6, 7).Value = 2.56
$ws.Cells.Item(6, 9).Value = 3.8
$ws.Cells.Item(6, 13).Value = 1.11
$ws.Cells.Item(6, 14).Value = 2.76
$ws.Cells.Item(6, 16).Value = 1.59
$ws.Cells.Item(6, 18).Value = 1.21
$ws.Cells.Item(6, 20).Value = 1.97
$ws.Cells.Item(6, 21).Value = 1.86
$ws.Cells.Item(6, 28).Value = 9.6
$ws.Cells.Item(6, 32).Value = 18
$ws.Cells.Item(6, 34).Value = 23
$ws.Cells.Item(6, 36).Value = 40
$ws.Cells.Item(6, 37).Value = 36
$ws.Cells.Item(6, 39).Value = 190
$ws.Cells.Item(6, 40).Value = 1000
$ws.Cells.Item(7, 6).Value = 2.58
$ws.Cells.Item(7, 7).Value = 2.84
$ws.Cells.Item(7, 8).Value = 3
$ws.Cells.Item(7, 9).Value = 3.4
$ws.Cells.Item(7, 10).Value = 3.1
$ws.Cells.Item(7, 14).Value = 2.8
$ws.Cells.Item(7, 20).Value = 1.98
$ws.Cells.Item(7, 21).Value = 1.88
$ws.Cells.Item(7, 22).Value = 1.42
$ws.Cells.Item(7, 23).Value = 1.54
$ws.Cells.Item(7, 26).Value = 29
$ws.Cells.Item(7, 27).Value = 85
$ws.Cells.Item(7, 30).Value = 20
$ws.Cells.Item(7, 31).Value = 65
$ws.Cells.Item(7, 34).Value = 29
$ws.Cells.Item(7, 35).Value = 95

# --- 2) Append new rows 8-11 ---

# Row 8
$ws.Cells.Item(8, 1).Value = "Brazilian Serie A"
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "2025-10-15"
$ws.Cells.Item(8, 2).Style = "Normal"
$ws.Cells.Item(8, 3).Value = "21:30:00"
$ws.Cells.Item(8, 4).Value = "Santos"
$ws.Cells.Item(8, 5).Value = "Corinthians"
$ws.Cells.Item(8, 6).Value = 2.22
$ws.Cells.Item(8, 7).Value = 2.34
$ws.Cells.Item(8, 8).Value = 3.85
$ws.Cells.Item(8, 9).Value = 4.3
$ws.Cells.Item(8, 10).Value = 3.2
$ws.Cells.Item(8, 11).Value = 3.3
$ws.Cells.Item(8, 12).Value = 1.6
$ws.Cells.Item(8, 13).Value = 1.13
$ws.Cells.Item(8, 14).Value = 2.58
$ws.Cells.Item(8, 15).Value = 1.55
$ws.Cells.Item(8, 16).Value = 1.53
$ws.Cells.Item(8, 17).Value = 2.62
$ws.Cells.Item(8, 18).Value = 1.19
$ws.Cells.Item(8, 19).Value = 5.7
$ws.Cells.Item(8, 20).Value = 2.14
$ws.Cells.Item(8, 21).Value = 1.77
$ws.Cells.Item(8, 22).Value = 1.25
$ws.Cells.Item(8, 23).Value = 1.58
$ws.Cells.Item(8, 24).Value = 8.6
$ws.Cells.Item(8, 25).Value = 11
$ws.Cells.Item(8, 26).Value = 980
$ws.Cells.Item(8, 27).Value = 120
$ws.Cells.Item(8, 28).Value = 7.2
$ws.Cells.Item(8, 29).Value = 7.4
$ws.Cells.Item(8, 30).Value = 980
$ws.Cells.Item(8, 31).Value = 70
$ws.Cells.Item(8, 32).Value = 13.5
$ws.Cells.Item(8, 33).Value = 980
$ws.Cells.Item(8, 34).Value = 980
$ws.Cells.Item(8, 35).Value = 110
$ws.Cells.Item(8, 36).Value = 980
$ws.Cells.Item(8, 37).Value = 980
$ws.Cells.Item(8, 38).Value = 65
$ws.Cells.Item(8, 39).Value = 240
$ws.Cells.Item(8, 40).Value = 980
$ws.Cells.Item(8, 41).Value = 130

# Row 9
$ws.Cells.Item(9, 1).Value = "Brazilian Serie A"
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "2025-10-15"
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 3).Value = "21:30:00"
$ws.Cells.Item(9, 4).Value = "Fortaleza EC"
$ws.Cells.Item(9, 5).Value = "Vasco Da Gama"
$ws.Cells.Item(9, 6).Value = 2.3
$ws.Cells.Item(9, 7).Value = 2.46
$ws.Cells.Item(9, 8).Value = 3.35
$ws.Cells.Item(9, 9).Value = 3.7
$ws.Cells.Item(9, 10).Value = 3.4
$ws.Cells.Item(9, 11).Value = 3.5
$ws.Cells.Item(9, 12).Value = 1.47
$ws.Cells.Item(9, 13).Value = 1.04
$ws.Cells.Item(9, 14).Value = 1.75
$ws.Cells.Item(9, 15).Value = 1.38
$ws.Cells.Item(9, 16).Value = 1.75
$ws.Cells.Item(9, 17).Value = 2.16
$ws.Cells.Item(9, 18).Value = 1.24
$ws.Cells.Item(9, 19).Value = 3.7
$ws.Cells.Item(9, 20).Value = 1.6
$ws.Cells.Item(9, 21).Value = 1.84
$ws.Cells.Item(9, 22).Value = 1.37
$ws.Cells.Item(9, 23).Value = 1.68
$ws.Cells.Item(9, 24).Value = 14
$ws.Cells.Item(9, 25).Value = 17
$ws.Cells.Item(9, 26).Value = 34
$ws.Cells.Item(9, 27).Value = 80
$ws.Cells.Item(9, 28).Value = 12.5
$ws.Cells.Item(9, 29).Value = 10.5
$ws.Cells.Item(9, 30).Value = 21
$ws.Cells.Item(9, 31).Value = 980
$ws.Cells.Item(9, 32).Value = 21
$ws.Cells.Item(9, 33).Value = 16.5
$ws.Cells.Item(9, 34).Value = 27
$ws.Cells.Item(9, 35).Value = 70
$ws.Cells.Item(9, 36).Value = 980
$ws.Cells.Item(9, 37).Value = 980
$ws.Cells.Item(9, 38).Value = 980
$ws.Cells.Item(9, 39).Value = 1000
$ws.Cells.Item(9, 40).Value = 1000
$ws.Cells.Item(9, 41).Value = 1000

# Row 10
$ws.Cells.Item(10, 1).Value = "Brazilian Serie A"
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "2025-10-15"
$ws.Cells.Item(10, 2).Style = "Normal"
$ws.Cells.Item(10, 3).Value = "21:30:00"
$ws.Cells.Item(10, 4).Value = "Atletico MG"
$ws.Cells.Item(10, 5).Value = "Cruzeiro MG"
$ws.Cells.Item(10, 6).Value = 3.05
$ws.Cells.Item(10, 7).Value = 3.35
$ws.Cells.Item(10, 8).Value = 2.6
$ws.Cells.Item(10, 9).Value = 2.76
$ws.Cells.Item(10, 10).Value = 3.05
$ws.Cells.Item(10, 11).Value = 3.25
$ws.Cells.Item(10, 12).Value = 1.57
$ws.Cells.Item(10, 13).Value = 1.12
$ws.Cells.Item(10, 14).Value = 2.68
$ws.Cells.Item(10, 15).Value = 1.52
$ws.Cells.Item(10, 16).Value = 1.57
$ws.Cells.Item(10, 17).Value = 2.58
$ws.Cells.Item(10, 18).Value = 1.21
$ws.Cells.Item(10, 19).Value = 5.1
$ws.Cells.Item(10, 20).Value = 2.04
$ws.Cells.Item(10, 21).Value = 1.82
$ws.Cells.Item(10, 22).Value = 1.56
$ws.Cells.Item(10, 23).Value = 1.45
$ws.Cells.Item(10, 24).Value = 9.199999999999999
$ws.Cells.Item(10, 25).Value = 8.4
$ws.Cells.Item(10, 26).Value = 980
$ws.Cells.Item(10, 27).Value = 980
$ws.Cells.Item(10, 28).Value = 9.4
$ws.Cells.Item(10, 29).Value = 7.2
$ws.Cells.Item(10, 30).Value = 13.5
$ws.Cells.Item(10, 31).Value = 980
$ws.Cells.Item(10, 32).Value = 980
$ws.Cells.Item(10, 33).Value = 980
$ws.Cells.Item(10, 34).Value = 980
$ws.Cells.Item(10, 35).Value = 65
$ws.Cells.Item(10, 36).Value = 65
$ws.Cells.Item(10, 37).Value = 980
$ws.Cells.Item(10, 38).Value = 70
$ws.Cells.Item(10, 39).Value = 210
$ws.Cells.Item(10, 40).Value = 60
$ws.Cells.Item(10, 41).Value = 980

# Row 11
$ws.Cells.Item(11, 1).Value = "US United Soccer League"
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "2025-10-15"
$ws.Cells.Item(11, 2).Style = "Normal"
$ws.Cells.Item(11, 3).Value = "23:00:00"
$ws.Cells.Item(11, 4).Value = "Orange County Blues"
$ws.Cells.Item(11, 5).Value = "San Antonio FC"
$ws.Cells.Item(11, 6).Value = 2.36
$ws.Cells.Item(11, 7).Value = 2.7
$ws.Cells.Item(11, 8).Value = 2.82
$ws.Cells.Item(11, 9).Value = 3.3
$ws.Cells.Item(11, 10).Value = 3.35
$ws.Cells.Item(11, 11).Value = 3.95
$ws.Cells.Item(11, 12).Value = 1.01
$ws.Cells.Item(11, 13).Value = 1.05
$ws.Cells.Item(11, 14).Value = 3.65
$ws.Cells.Item(11, 15).Value = 1.3
$ws.Cells.Item(11, 16).Value = 1.92
$ws.Cells.Item(11, 17).Value = 1.87
$ws.Cells.Item(11, 18).Value = 1.35
$ws.Cells.Item(11, 19).Value = 3.2
$ws.Cells.Item(11, 20).Value = 1.7
$ws.Cells.Item(11, 21).Value = 2.12
$ws.Cells.Item(11, 22).Value = 1.43
$ws.Cells.Item(11, 23).Value = 1.58
$ws.Cells.Item(11, 24).Value = 18.5
$ws.Cells.Item(11, 25).Value = 13
$ws.Cells.Item(11, 26).Value = 980
$ws.Cells.Item(11, 27).Value = 55
$ws.Cells.Item(11, 28).Value = 11.5
$ws.Cells.Item(11, 29).Value = 8.6
$ws.Cells.Item(11, 30).Value = 14
$ws.Cells.Item(11, 31).Value = 980
$ws.Cells.Item(11, 32).Value = 18
$ws.Cells.Item(11, 33).Value = 12.5
$ws.Cells.Item(11, 34).Value = 18
$ws.Cells.Item(11, 35).Value = 980
$ws.Cells.Item(11, 36).Value = 980
$ws.Cells.Item(11, 37).Value = 980
$ws.Cells.Item(11, 38).Value = 980
$ws.Cells.Item(11, 39).Value = 1000
$ws.Cells.Item(11, 40).Value = 26
$ws.Cells.Item(11, 41).Value = 36

Write-Host "Update complete"